# Update timing/memory metric values on the active sheet to reflect the
# latest benchmark run (background fix + migration to new template).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.000347365
$ws.Range("F2").Value = 0.025842378
$ws.Range("G2").Value = 0.0006340764801227152

$ws.Range("E3").Value = 0.006656237
$ws.Range("F3").Value = 0.011940453
$ws.Range("G3").Value = 0.007307757475146198
